$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "578.18") are written as text, matching the original inlineStr cells,
# then restore the default style so no stray formatting is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.352.39"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "3.227.45"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "578.18"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").Value = "183.76"
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "3.223.93"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("E11").Value = "  -2.34%  "
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "3.780.46"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "27.66"
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").Value = "67.433.46"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").Value = "3.225.12"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").Value = "13.43"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").Value = "395.26"
$ws.Range("E21").Value = "  +3.56%  "
$ws.Range("D22").Value = "7.54"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "71.09"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").Value = "0.515"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -2.74%  "
$ws.Range("D27").Value = "0.185"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("D28").Value = "9.50"
$ws.Range("E28").Value = "  -3.29%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").Value = "5.54"
$ws.Range("E31").Value = "  -5.09%  "
$ws.Range("D32").Value = "22.55"
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("D33").Value = "6.96"
$ws.Range("E33").Value = "  -3.29%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.25"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "160.51"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").Value = "26.38"
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "0.802"
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").Value = "6.51"
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("D43").Value = "2.46"
$ws.Range("E43").Value = "  -5.91%  "
$ws.Range("D44").Value = "0.0683"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").Value = "40.42"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("D46").Value = "2.595.78"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").Value = "334.35"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("D48").Value = "24.54"
$ws.Range("E48").Value = "  -3.42%  "
$ws.Range("E49").Value = "  -2.70%  "
$ws.Range("D50").Value = "6.26"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  -1.06%  "

$priceRange.Style = "Normal"
